{"js": "const pairs = [\n  [\"2023-09-15 Friday\", \"2023-09-16 Saturday\"],\n  [\"31\u00f79=3, 4\", \"53\u00f74=13, 1\"],\n  [\"82\u00f74=20, 2\", \"31\u00f77=4, 3\"],\n  [\"19\u00f79=2, 1\", \"11\u00f74=2, 3\"],\n  [\"56\u00f78=7, 0\", \"85\u00f72=42, 1\"],\n  [\"93\u00f79=10, 3\", \"82\u00f79=9, 1\"],\n  [\"64\u00f75=12, 4\", \"89\u00f72=44, 1\"],\n  [\"75\u00f76=12, 3\", \"74\u00f74=18, 2\"],\n  [\"40\u00f75=8, 0\", \"94\u00f72=47, 0\"],\n  [\"23\u00f76=3, 5\", \"61\u00f78=7, 5\"],\n  [\"49\u00f73=16, 1\", \"96\u00f79=10, 6\"],\n  [\"35\u00f73=11, 2\", \"76\u00f75=15, 1\"],\n  [\"42\u00f75=8, 2\", \"42\u00f76=7, 0\"],\n  [\"28\u00f79=3, 1\", \"88\u00f72=44, 0\"],\n  [\"62\u00f72=31, 0\", \"84\u00f78=10, 4\"],\n  [\"90\u00f76=15, 0\", \"72\u00f74=18, 0\"],\n  [\"92\u00f76=15, 2\", \"40\u00f75=8, 0\"],\n  [\"42\u00f74=10, 2\", \"65\u00f78=8, 1\"],\n  [\"28\u00f77=4, 0\", \"97\u00f79=10, 7\"],\n  [\"81\u00f74=20, 1\", \"15\u00f73=5, 0\"],\n  [\"44\u00f77=6, 2\", \"73\u00f78=9, 1\"],\n  [\"93\u00f76=15, 3\", \"21\u00f73=7, 0\"],\n  [\"38\u00f76=6, 2\", \"96\u00f77=13, 5\"],\n  [\"85\u00f78=10, 5\", \"13\u00f79=1, 4\"],\n  [\"22\u00f72=11, 0\", \"59\u00f72=29, 1\"],\n  [\"25\u00f79=2, 7\", \"35\u00f74=8, 3\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-09-15 Friday\", \"2023-09-16 Saturday\"),\n    @(\"31\u00f79=3, 4\", \"53\u00f74=13, 1\"),\n    @(\"82\u00f74=20, 2\", \"31\u00f77=4, 3\"),\n    @(\"19\u00f79=2, 1\", \"11\u00f74=2, 3\"),\n    @(\"56\u00f78=7, 0\", \"85\u00f72=42, 1\"),\n    @(\"93\u00f79=10, 3\", \"82\u00f79=9, 1\"),\n    @(\"64\u00f75=12, 4\", \"89\u00f72=44, 1\"),\n    @(\"75\u00f76=12, 3\", \"74\u00f74=18, 2\"),\n    @(\"40\u00f75=8, 0\", \"94\u00f72=47, 0\"),\n    @(\"23\u00f76=3, 5\", \"61\u00f78=7, 5\"),\n    @(\"49\u00f73=16, 1\", \"96\u00f79=10, 6\"),\n    @(\"35\u00f73=11, 2\", \"76\u00f75=15, 1\"),\n    @(\"42\u00f75=8, 2\", \"42\u00f76=7, 0\"),\n    @(\"28\u00f79=3, 1\", \"88\u00f72=44, 0\"),\n    @(\"62\u00f72=31, 0\", \"84\u00f78=10, 4\"),\n    @(\"90\u00f76=15, 0\", \"72\u00f74=18, 0\"),\n    @(\"92\u00f76=15, 2\", \"40\u00f75=8, 0\"),\n    @(\"42\u00f74=10, 2\", \"65\u00f78=8, 1\"),\n    @(\"28\u00f77=4, 0\", \"97\u00f79=10, 7\"),\n    @(\"81\u00f74=20, 1\", \"15\u00f73=5, 0\"),\n    @(\"44\u00f77=6, 2\", \"73\u00f78=9, 1\"),\n    @(\"93\u00f76=15, 3\", \"21\u00f73=7, 0\"),\n    @(\"38\u00f76=6, 2\", \"96\u00f77=13, 5\"),\n    @(\"85\u00f78=10, 5\", \"13\u00f79=1, 4\"),\n    @(\"22\u00f72=11, 0\", \"59\u00f72=29, 1\"),\n    @(\"25\u00f79=2, 7\", \"35\u00f74=8, 3\"),\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}"}
